$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.923.18"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "2.400.42"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.15%  "

$cell = $ws.Range("D5")
$cell.Value = "'554.86"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

$cell = $ws.Range("D6")
$cell.Value = "'134.77"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.82%  "

$ws.Range("E7").Value = "  +0.15%  "

$cell = $ws.Range("D8")
$cell.Value = "'0.585"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.38%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.105"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.65%  "

$cell = $ws.Range("D10")
$cell.Value = "'5.58"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("E11").Value = "  +0.25%  "

$cell = $ws.Range("D12")
$cell.Value = "'0.346"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.41%  "

$cell = $ws.Range("D13")
$cell.Value = "'24.63"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("D14").Value = "2.834.17"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "59.944.29"
$ws.Range("E15").Value = "  +0.39%  "

$cell = $ws.Range("D16")
$cell.Value = "'0.0000137"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "2.419.38"
$ws.Range("E17").Value = "  +0.23%  "

$cell = $ws.Range("D18")
$cell.Value = "'11.16"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "

$cell = $ws.Range("D19")
$cell.Value = "'4.49"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.60%  "

$cell = $ws.Range("D20")
$cell.Value = "'325.92"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "

$cell = $ws.Range("D21")
$cell.Value = "'6.74"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "

$cell = $ws.Range("D22")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$cell = $ws.Range("D23")
$cell.Value = "'64.69"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.05%  "

$cell = $ws.Range("D24")
$cell.Value = "'0.172"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "

$cell = $ws.Range("D25")
$cell.Value = "'8.49"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.56%  "

$ws.Range("E26").Value = "  +0.03%  "

$cell = $ws.Range("D27")
$cell.Value = "'1.38"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "

$cell = $ws.Range("D28")
$cell.Value = "'1.79"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").Value = "0.0₃0763"
$ws.Range("E29").Value = "  -1.43%  "

$cell = $ws.Range("D30")
$cell.Value = "'170.68"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.38%  "

$cell = $ws.Range("D31")
$cell.Value = "'6.09"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("E32").Value = "  +6.23%  "

$cell = $ws.Range("D33")
$cell.Value = "'0.401"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.10%  "

$cell = $ws.Range("D34")
$cell.Value = "'18.34"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D35")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D36")
$cell.Value = "'1.33"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.15%  "

$ws.Range("E37").Value = "  +0.10%  "

$cell = $ws.Range("D38")
$cell.Value = "'4.15"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$cell = $ws.Range("D39")
$cell.Value = "'322.87"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

$cell = $ws.Range("D40")
$cell.Value = "'1.59"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "

$cell = $ws.Range("D41")
$cell.Value = "'38.50"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.75%  "

$cell = $ws.Range("D42")
$cell.Value = "'147.80"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +6.56%  "

$cell = $ws.Range("D43")
$cell.Value = "'3.54"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -3.57%  "

$cell = $ws.Range("D44")
$cell.Value = "'0.0965"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "

$cell = $ws.Range("D45")
$cell.Value = "'19.66"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.23%  "

$cell = $ws.Range("D46")
$cell.Value = "'0.0514"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "

$cell = $ws.Range("D47")
$cell.Value = "'0.574"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "

$cell = $ws.Range("D48")
$cell.Value = "'0.0220"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.74%  "

$cell = $ws.Range("D49")
$cell.Value = "'11.04"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.19%  "

$cell = $ws.Range("D50")
$cell.Value = "'1.56"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("E51").Value = "  -0.87%  "
